# Auto-generated: update market-price derived columns (H-N) across profession sheets
# per scheduled market data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1818.0571
$ws.Range("I70").Value = 1610.3334
$ws.Range("J70").Value = 2129.6428
$ws.Range("K70").Value = 4831.0002
$ws.Range("L70").Value = 6388.928400000001
$ws.Range("M70").Value = -4561.0002
$ws.Range("N70").Value = -6928.928400000001
$ws.Range("H73").Value = 1818.0571
$ws.Range("I73").Value = 1610.3334
$ws.Range("J73").Value = 2129.6428
$ws.Range("K73").Value = 4831.0002
$ws.Range("L73").Value = 6388.928400000001
$ws.Range("M73").Value = -3895.0002
$ws.Range("N73").Value = -8260.928400000001
$ws.Range("H76").Value = 4357.143
$ws.Range("I76").Value = 3441.1765
$ws.Range("K76").Value = 3441.1765
$ws.Range("M76").Value = -3126.1765
$ws.Range("H79").Value = 4357.143
$ws.Range("I79").Value = 3441.1765
$ws.Range("K79").Value = 3441.1765
$ws.Range("M79").Value = -2349.1765
$ws.Range("H82").Value = 247.33333
$ws.Range("I82").Value = 247.33333
$ws.Range("K82").Value = 741.99999
$ws.Range("M82").Value = -335.99999
$ws.Range("H85").Value = 247.33333
$ws.Range("I85").Value = 247.33333
$ws.Range("K85").Value = 741.99999
$ws.Range("M85").Value = 662.00001
$ws.Range("H86").Value = 31198.354
$ws.Range("I86").Value = 45188.87
$ws.Range("J86").Value = 1945.4546
$ws.Range("K86").Value = 45188.87
$ws.Range("L86").Value = 1945.4546
$ws.Range("M86").Value = -44065.87
$ws.Range("N86").Value = -4191.4546
$ws.Range("H88").Value = 2936.625
$ws.Range("I88").Value = 1831
$ws.Range("J88").Value = 3600
$ws.Range("K88").Value = 1831
$ws.Range("L88").Value = 3600
$ws.Range("M88").Value = -1425
$ws.Range("N88").Value = -4412
$ws.Range("H89").Value = 31198.354
$ws.Range("I89").Value = 45188.87
$ws.Range("J89").Value = 1945.4546
$ws.Range("K89").Value = 225944.35
$ws.Range("L89").Value = 9727.273000000001
$ws.Range("M89").Value = -220328.35
$ws.Range("N89").Value = -20959.273
$ws.Range("H91").Value = 2936.625
$ws.Range("I91").Value = 1831
$ws.Range("J91").Value = 3600
$ws.Range("K91").Value = 1831
$ws.Range("L91").Value = 3600
$ws.Range("M91").Value = -427
$ws.Range("N91").Value = -6408
$ws.Range("H92").Value = 331
$ws.Range("I92").Value = 355.73334
$ws.Range("K92").Value = 355.73334
$ws.Range("M92").Value = 892.26666
$ws.Range("H97").Value = 5500
$ws.Range("J97").Value = 5500
$ws.Range("L97").Value = 16500
$ws.Range("N97").Value = -17492
$ws.Range("H100").Value = 2250.1738
$ws.Range("I100").Value = 1608.8572
$ws.Range("J100").Value = 3247.7778
$ws.Range("K100").Value = 1608.8572
$ws.Range("L100").Value = 3247.7778
$ws.Range("M100").Value = -1067.8572
$ws.Range("N100").Value = -4329.7778
$ws.Range("H103").Value = 1500
$ws.Range("J103").Value = 1500
$ws.Range("L103").Value = 4500
$ws.Range("N103").Value = -5672
$ws.Range("H105").Value = 32666.666
$ws.Range("J105").Value = 32666.666
$ws.Range("L105").Value = 32666.666
$ws.Range("N105").Value = -39654.666
$ws.Range("H135").Value = 125002410
$ws.Range("I135").Value = 55557100
$ws.Range("J135").Value = 333338340
$ws.Range("K135").Value = 500013900
$ws.Range("L135").Value = 3000045060
$ws.Range("M135").Value = -500011365
$ws.Range("N135").Value = -3000050130
$ws.Range("H138").Value = 3373.0293
$ws.Range("J138").Value = 4080.5908
$ws.Range("L138").Value = 12241.7724
$ws.Range("N138").Value = -22521.7724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 507543.7
$ws.Range("I32").Value = 22686.469
$ws.Range("J32").Value = 2335082.5
$ws.Range("K32").Value = 22686.469
$ws.Range("L32").Value = 2335082.5
$ws.Range("M32").Value = -22399.469
$ws.Range("N32").Value = -2335656.5
$ws.Range("H69").Value = 76000
$ws.Range("J69").Value = 76000
$ws.Range("L69").Value = 76000
$ws.Range("N69").Value = -77498
$ws.Range("H72").Value = 76000
$ws.Range("J72").Value = 76000
$ws.Range("L72").Value = 228000
$ws.Range("N72").Value = -235488
$ws.Range("H74").Value = 4749.2256
$ws.Range("I74").Value = 1649.9231
$ws.Range("J74").Value = 20865.6
$ws.Range("K74").Value = 1649.9231
$ws.Range("L74").Value = 20865.6
$ws.Range("M74").Value = -775.9231
$ws.Range("N74").Value = -22613.6
$ws.Range("H77").Value = 4749.2256
$ws.Range("I77").Value = 1649.9231
$ws.Range("J77").Value = 20865.6
$ws.Range("K77").Value = 8249.6155
$ws.Range("L77").Value = 104328
$ws.Range("M77").Value = -3881.6155
$ws.Range("N77").Value = -113064
$ws.Range("H122").Value = 9616853
$ws.Range("I122").Value = 1443.2222
$ws.Range("J122").Value = 31251524
$ws.Range("K122").Value = 4329.6666
$ws.Range("L122").Value = 93754572
$ws.Range("M122").Value = -1879.6666
$ws.Range("N122").Value = -93759472

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 353.5
$ws.Range("I16").Value = 353.5
$ws.Range("K16").Value = 353.5
$ws.Range("M16").Value = -183.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 56542
$ws.Range("J18").Value = 56542
$ws.Range("L18").Value = 56542
$ws.Range("N18").Value = -57002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 598.2
$ws.Range("I107").Value = 275.02084
$ws.Range("J107").Value = 2814.2856
$ws.Range("K107").Value = 825.0625200000001
$ws.Range("L107").Value = 8442.856800000001
$ws.Range("M107").Value = 1094.93748
$ws.Range("N107").Value = -12282.8568
$ws.Range("H113").Value = 340.12
$ws.Range("I113").Value = 368.2258
$ws.Range("J113").Value = 327.49277
$ws.Range("K113").Value = 1104.6774
$ws.Range("L113").Value = 982.47831
$ws.Range("M113").Value = 1065.3226
$ws.Range("N113").Value = -5322.47831
$ws.Range("H131").Value = 722.7805
$ws.Range("I131").Value = 212.3125
$ws.Range("J131").Value = 1049.48
$ws.Range("K131").Value = 636.9375
$ws.Range("L131").Value = 3148.44
$ws.Range("M131").Value = 4403.0625
$ws.Range("N131").Value = -13228.44

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2707.25
$ws.Range("J126").Value = 3397.1538
$ws.Range("L126").Value = 10191.4614
$ws.Range("N126").Value = -15131.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 20133.285
$ws.Range("I61").Value = 22852.047
$ws.Range("J61").Value = 11977
$ws.Range("K61").Value = 22852.047
$ws.Range("L61").Value = 11977
$ws.Range("M61").Value = -22650.047
$ws.Range("N61").Value = -12381
$ws.Range("H113").Value = 20133.285
$ws.Range("I113").Value = 22852.047
$ws.Range("J113").Value = 11977
$ws.Range("K113").Value = 22852.047
$ws.Range("L113").Value = 11977
$ws.Range("M113").Value = -20682.047
$ws.Range("N113").Value = -16317
$ws.Range("H132").Value = 4947.0625
$ws.Range("J132").Value = 8195.571
$ws.Range("L132").Value = 24586.713
$ws.Range("N132").Value = -29646.713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 19087.5
$ws.Range("J86").Value = 19087.5
$ws.Range("L86").Value = 19087.5
$ws.Range("N86").Value = -21333.5
$ws.Range("H89").Value = 19087.5
$ws.Range("J89").Value = 19087.5
$ws.Range("L89").Value = 95437.5
$ws.Range("N89").Value = -106669.5
